$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 189 (this shifts the existing rows 189-218 down to 190-219,
# carrying formatting, such as the date style on column D, along with them).
$ws.Rows.Item(189).Insert()

# Copy the static / repeated columns from the row right below (which is the
# former row 189, now shifted to row 190) since every data row in this sheet
# shares the same Mercado/Region/Categoria/etc. values.
$cols = @(1, 2, 3, 5, 6, 7, 8, 9, 14, 17, 18)
foreach ($col in $cols) {
    $ws.Cells.Item(189, $col).Value = $ws.Cells.Item(190, $col).Value()
}

# New row-specific data for row 189.
$ws.Cells.Item(189, 4).Value = 44644   # D189 Fecha
$ws.Cells.Item(189, 10).Value = 65     # J189 Volumen
$ws.Cells.Item(189, 11).Value = 8000   # K189 Precio minimo
$ws.Cells.Item(189, 12).Value = 8000   # L189 Precio maximo
$ws.Cells.Item(189, 13).Value = 8000   # M189 Precio promedio ponderado
$ws.Cells.Item(189, 15).Value = "Región de La Araucanía"  # O189 Origen
$ws.Cells.Item(189, 16).Value = 8000   # P189 Precio $/Kg
